# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New computed "K" values (column G) replacing the old Strike# derived values.
$newK = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 2
    9  = 1
    10 = 2
    11 = 0
    12 = 1
    13 = 4
    14 = 3
    15 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
